$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.532.42'
$ws.Range('D3').Value = '1.626.98'
$ws.Range('E3').Value = '  -0.54%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '212.03'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.18%  '
$ws.Range('E6').Value = '  -0.36%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.18%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '23.36'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +1.93%  '
$ws.Range('E9').Value = '  +2.25%  '
$ws.Range('E10').Value = '  +0.23%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0877'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.91%  '
$ws.Range('D12').Value = '1.857.14'
$ws.Range('E12').Value = '  -0.57%  '
$ws.Range('D13').Value = '1.631.52'
$ws.Range('E13').Value = '  -0.23%  '
$ws.Range('E14').Value = '  +0.34%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.553'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -1.22%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '65.38'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +1.35%  '
$ws.Range('D17').Value = '27.512.19'
$ws.Range('E17').Value = '  -0.36%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '229.98'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.34%  '
$ws.Range('E19').Value = '  -0.54%  '
$ws.Range('E20').Value = '  -2.21%  '
$ws.Range('E21').Value = '  -0.07%  '
$ws.Range('E22').Value = '  +4.64%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.35'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.78%  '
$ws.Range('E24').Value = '  +8.77%  '
$ws.Range('E25').Value = '  -0.76%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '6.88'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.52%  '
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('B28').Value = 'BinanceUSD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.00'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.13%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.52'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.64%  '
$ws.Range('E30').Value = '  -0.29%  '
$ws.Range('E31').Value = '  -0.25%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.28'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -0.58%  '
$ws.Range('D33').Value = '1.468.19'
$ws.Range('E33').Value = '  +0.98%  '
$ws.Range('E34').Value = '  -1.65%  '
$ws.Range('E35').Value = '  -0.91%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.941'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +5.51%  '
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.876'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.12%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0167'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.54%  '
$ws.Range('E40').Value = '  -1.62%  '
$ws.Range('E41').Value = '  +2.12%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.999'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.16%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '67.77'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -2.77%  '
$ws.Range('E44').Value = '  +0.41%  '
$ws.Range('E45').Value = '  -1.61%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '5.35'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -4.31%  '
$ws.Range('E47').Value = '  +3.72%  '
$ws.Range('D48').Value = '1.767.05'
$ws.Range('E48').Value = '  -0.65%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '87.35'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +1.34%  '
$ws.Range('E50').Value = '  -0.87%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0994'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.98%  '
